$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ', '

    if ($parts.Length -lt 2) { continue }

    if ($parts[0].ToLower() -ne 'system') { continue }

    $first = $parts[0]
    $last = $parts[$parts.Length - 1]
    $parts[0] = $last
    $parts[$parts.Length - 1] = $first

    $newVal = $parts -join ', '
    $cell.Value2 = $newVal
}
